$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: yChannel changes from POS/RelOriginAlt -> CTUN/Alt
$ws.Range("K12").Value = "CTUN/Alt"

# Row 13: yChannel changes from CTUN/DSAlt -> CTUN/SAlt, LabelOverride changes from "Des Alt" -> "Sonar Altitude"
$ws.Range("K13").Value = "CTUN/SAlt"
$ws.Range("R13").Value = "Sonar Altitude"

# New row 14: duplicate of the old row 13 (Des Alt plot), with new label/channel
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = "S"
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("E14").Value = "Time [ s ]"
$ws.Range("F14").Value = "Altitude~[~m~]"
$ws.Range("G14").Value = "Vertical"
$ws.Range("K14").Value = "CTUN/DSAlt"
$ws.Range("O14").Value = 1
$ws.Range("R14").Value = "Des Altitude"
$ws.Range("S14").Value = "m"

# Update sheet view: scroll position and selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("R14").Select()
